$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.790.14'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '1.886.73'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '239.05'
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '0.4755'
$ws.Range("E7").Value = '  +1.53%  '
$ws.Range("E8").Value = '  +5.09%  '
$ws.Range("D9").Value = '0.06580'
$ws.Range("E9").Value = '  +4.11%  '
$ws.Range("D10").Value = '18.72'
$ws.Range("E10").Value = '  +9.08%  '
$ws.Range("D11").Value = '99.37'
$ws.Range("E11").Value = '  +17.99%  '
$ws.Range("D12").Value = '1.886.90'
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").Value = '0.07595'
$ws.Range("D14").Value = '5.131'
$ws.Range("E14").Value = '  +3.72%  '
$ws.Range("D15").Value = '0.6623'
$ws.Range("E15").Value = '  +5.49%  '
$ws.Range("D16").Value = '307.89'
$ws.Range("E16").Value = '  +33.91%  '
$ws.Range("D17").Value = '30.804.50'
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").Value = '13.18'
$ws.Range("E18").Value = '  +5.04%  '
$ws.Range("D19").Value = '1.0000'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '0.000007579'
$ws.Range("E20").Value = '  +3.55%  '
$ws.Range("D21").Value = '2.135.90'
$ws.Range("E21").Value = '  +2.20%  '
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '5.125'
$ws.Range("E23").Value = '  +3.78%  '
$ws.Range("D24").Value = '6.218'
$ws.Range("E24").Value = '  +5.22%  '
$ws.Range("D25").Value = '9.304'
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").Value = '167.45'
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("D27").Value = '20.52'
$ws.Range("E27").Value = '  +14.35%  '
$ws.Range("D28").Value = '1.948'
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("E29").Value = '  +6.15%  '
$ws.Range("D30").Value = '1.355'
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("D31").Value = '4.182'
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("D32").Value = '3.983'
$ws.Range("E32").Value = '  +3.88%  '
$ws.Range("D33").Value = '0.05086'
$ws.Range("E33").Value = '  +3.98%  '
$ws.Range("D34").Value = '1.171'
$ws.Range("E34").Value = '  +2.59%  '
$ws.Range("D35").Value = '0.7293'
$ws.Range("E35").Value = '  +3.51%  '
$ws.Range("D36").Value = '2.716'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").Value = '0.01957'
$ws.Range("E37").Value = '  +2.40%  '
$ws.Range("D38").Value = '2.704'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("D39").Value = '2.074'
$ws.Range("E39").Value = '  +6.19%  '
$ws.Range("D40").Value = '0.9083'
$ws.Range("E40").Value = '  +4.08%  '
$ws.Range("D41").Value = '108.19'
$ws.Range("E41").Value = '  +2.30%  '
$ws.Range("D42").Value = '0.9997'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '0.4209'
$ws.Range("E43").Value = '  +3.67%  '
$ws.Range("D44").Value = '5.627'
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '7.387'
$ws.Range("E45").Value = '  +3.32%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '65.94'
$ws.Range("E46").Value = '  +7.09%  '
$ws.Range("D47").Value = '0.1230'
$ws.Range("E47").Value = '  +0.81%  '
$ws.Range("D48").Value = '8.986'
$ws.Range("E48").Value = '  +4.81%  '
$ws.Range("D49").Value = '34.84'
$ws.Range("E49").Value = '  +4.13%  '
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("D51").Value = '1.392'
$ws.Range("E51").Value = '  +2.05%  '
